$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not numeric) storage for price cells whose new value would
# otherwise be auto-parsed by Excel as a number, so they keep matching the
# original inline-string / text representation of the price column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values from the refreshed crypto price feed.
$ws.Range("D2").Value = "34.155.01"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.793.30"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "227.37"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  +3.19%  "
$ws.Range("D10").Value = "0.0693"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "2.050.97"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "11.54"
$ws.Range("E13").Value = "  +5.19%  "
$ws.Range("D14").Value = "1.796.00"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "34.122.04"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "68.00"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "245.15"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "4.12"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "161.84"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").Value = "16.33"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").Value = "1.443.62"
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("D36").Value = "0.650"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("E38").Value = "  +8.88%  "
$ws.Range("D39").Value = "1.03"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").Value = "80.88"
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("D41").Value = "0.930"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").Value = "2.36"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "2.70"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "13.33"
$ws.Range("E44").Value = "  +6.94%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "6.07"
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").Value = "0.0509"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0139"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").Value = "108.09"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").Value = "1.952.27"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").Value = "  -0.04%  "
